# Regenerate Report for Handback:
# bump the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the first data row (54a3c540-...) across the Overview, zh-cn and
# de-de sheets, simulating a re-run of the handback status report a
# little over a minute later.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the first file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-27 15:02:48"

# zh-cn sheet: Correspond Handoff Datetime (H) / Correspond Handback DateTime (K).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-27 15:02:44"
$wsZhCn.Range("K2").Value = "2016-08-27 15:03:06"

# de-de sheet: Correspond Handoff Datetime (H) / Correspond Handback DateTime (K).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-27 15:02:48"
$wsDeDe.Range("K2").Value = "2016-08-27 15:03:13"
